$d = $word.ActiveDocument

# Locate the exact phrase that changes ("vamos a continuar" -> "se procederá a
# continuar"). The surrounding text ("En esta sección " ... " con el
# desarrollo de la aplicación de tipo contador: ") stays the same, but in the
# target document that single run is split into three runs (identical
# formatting) around the replaced phrase. Toggling a character-formatting
# property while the replacement text is inserted forces the run to be split
# at the selection boundaries even though the property is reset back
# afterwards, matching the three-run shape shown in the diff.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("vamos a continuar", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Font.Bold = 1
    $rng.Text = "se procederá a continuar"
    $rng.Font.Bold = 0
}
